$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referenten")

# Delete the duplicate row 7 (A7/B7 duplicated A6/B6); remaining rows shift up.
$ws.Rows.Item(7).Delete()

# Make "Referenten" the active sheet, reset the frozen-pane scroll position to
# the top and select cell A7, matching the saved view state.
$ws.Activate()
$ws.Range("A2").Select()
$ws.Range("A7").Select()
